# Generate Report for Handoff
#
# Updates the "latest handoff" timestamps for the row corresponding to
# 7aa7752c-c2c4-4003-92d4-37f0044a715b.md, across the Overview sheet and
# both per-locale sheets (zh-cn, de-de), as if a fresh handoff report had
# just been generated for that file.

$wb = $excel.ActiveWorkbook

# Overview sheet: row 7 is the 7aa7752c... file; column G is
# "Latest HO Xliff Generate Date".
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G7").Value = "2016-09-04 12:45:49"

# zh-cn sheet: row 7 is the 7aa7752c... file; column H is
# "Latest Handoff Datetime".
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H7").Value = "2016-09-04 12:45:44"

# de-de sheet: row 7 is the 7aa7752c... file; column H is
# "Latest Handoff Datetime".
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H7").Value = "2016-09-04 12:45:49"
